$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "theorielessen applicaties"
$ws.Range("D17").Value = "indienen verslag + ppt"
